# Updates cryptos list values (Price / Volume(1h) columns, plus a few
# Coin/Link row swaps) to match the Sat Dec 16 08:52:01 UTC 2023 GitHub
# Actions refresh.
#
# Many "Price" column values are strings that look numeric (e.g. "6.14").
# Assigning such a string straight to Range.Value lets Excel auto-convert
# it to a real number (losing formatting like trailing zeros, e.g. "62.20"
# would become 62.2). To keep them as literal text - matching the source
# workbook, where every cell in these columns is a string - we prefix
# those values with a leading apostrophe (Excel's standard "force text /
# quote-prefix" convention), only where needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.157.39'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '2.241.63'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").Value = '''0.629'
$ws.Range("E6").Value = '  -2.20%  '
$ws.Range("D7").Value = '''74.25'
$ws.Range("E7").Value = '  -3.18%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '''0.618'
$ws.Range("E9").Value = '  -5.52%  '
$ws.Range("D10").Value = '''42.13'
$ws.Range("E10").Value = '  +4.98%  '
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("D12").Value = '''7.17'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '''14.43'
$ws.Range("E14").Value = '  -3.69%  '
$ws.Range("D15").Value = '''0.849'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '2.240.77'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '42.057.78'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '0.0₃0982'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''6.14'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '''72.11'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").Value = '''2.21'
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''230.03'
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").Value = '''8.86'
$ws.Range("E23").Value = '  +38.49%  '
$ws.Range("D25").Value = '''11.49'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  -5.34%  '
$ws.Range("D27").Value = '''2.29'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '''168.99'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '''20.65'
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("D31").Value = '''0.0818'
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("D32").Value = '''31.05'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("E35").Value = '  +10.00%  '
$ws.Range("D36").Value = '''4.47'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("E37").Value = '  +3.00%  '
$ws.Range("D38").Value = '''13.71'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -3.36%  '
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").Value = '''0.206'
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").Value = '''62.20'
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("D43").Value = '''106.69'
$ws.Range("E43").Value = '  -4.06%  '
$ws.Range("D44").Value = '''0.103'
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = '''0.996'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").Value = '''1.17'
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = '''2.29'
$ws.Range("E49").Value = '  +2.24%  '
$ws.Range("D50").Value = '''4.19'
$ws.Range("E50").Value = '  -6.91%  '
$ws.Range("E51").Value = '  +0.44%  '
